$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename TCID column (A2:A22) from TestCase_F1..TestCase_F21 to Notifications001..Notifications021
for ($r = 2; $r -le 22; $r++) {
    $n = $r - 1
    $label = "Notifications{0:D3}" -f $n
    $ws.Cells.Item($r, 1).Value = $label
}

# Runmode column (D2:D22): flip any "N" to "Y" so the whole column reads Y
for ($r = 2; $r -le 22; $r++) {
    $ws.Cells.Item($r, 4).Value = "Y"
}

# Widen column A slightly
$ws.Columns.Item(1).ColumnWidth = 16.16

# Move the view: scroll back to A1 and park the selection on C12
$ws.Activate()
$ws.Range("C12").Select()
